$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("O2").Value = 0.003189166373186978
$ws.Range("P2").Value = 0.003189166373186979
$ws.Range("S2").Value = 0.003189166373186978
$ws.Range("T2").Value = 0.003189166373186979

# Row 3
$ws.Range("M3").Value = 0.9980816666666668
$ws.Range("O3").Value = 0.05003796658465947
$ws.Range("P3").Value = 0.05003796658465948
$ws.Range("Q3").Value = 0.3242667526833334
$ws.Range("S3").Value = 0.05003796658465947
$ws.Range("T3").Value = 0.05003796658465948

# Row 4
$ws.Range("M4").Value = 15.66265066666667
$ws.Range("N4").Value = 46.987952
$ws.Range("O4").Value = 0.7852335303415662
$ws.Range("P4").Value = 0.7852335303415663
$ws.Range("Q4").Value = 5.088638575093333
$ws.Range("R4").Value = 45.79774717584
$ws.Range("S4").Value = 0.7852335303415662
$ws.Range("T4").Value = 0.7852335303415663

# Row 5
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.03576666666666667
$ws.Range("N5").Value = 0.1073
$ws.Range("O5").Value = 0.001793131094661246
$ws.Range("P5").Value = 0.001793131094661246
$ws.Range("Q5").Value = 0.01162023233333333
$ws.Range("R5").Value = 0.104582091
$ws.Range("S5").Value = 0.001793131094661246
$ws.Range("T5").Value = 0.001793131094661246

# Row 6
$ws.Range("M6").Value = 3.186375666666667
$ws.Range("N6").Value = 9.559127
$ws.Range("O6").Value = 0.1597462056059261
$ws.Range("P6").Value = 0.1597462056059261
$ws.Range("Q6").Value = 1.035221590343333
$ws.Range("R6").Value = 9.316994313090001
$ws.Range("S6").Value = 0.1597462056059261
$ws.Range("T6").Value = 0.1597462056059261
